$d = $word.ActiveDocument

function Find-ParaIndex($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# --- 1. Shorten the opening "Thank you" paragraph ---
$ok1 = $d.Content.Find.Execute(
    "Thank you for sending review reports for our manuscript. The reviewer’s comments were again useful and have improved the paper.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Thank you for sending review reports for our manuscript.",
    2)
if (-not $ok1) { Write-Output "WARN: step 1 find/replace did not match" }

# --- 2. Replace the "We enclose a revised version..." paragraph text, then
#        insert two new paragraphs after it (reworded "enclose" sentence +
#        new "We hope..." sentence). ---
$ok2 = $d.Content.Find.Execute(
    "We enclose a revised version of the manuscript (clean and with tracked changes), which addresses the reviewer’s comments and suggestions. The reviewer’s comments are shown below with our response to each comment in italics.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We have revised the paper in light of the reviewers’ comments. In particular, we have revised the tables for clarity, and improved their captioning.",
    2)
if (-not $ok2) { Write-Output "WARN: step 2 find/replace did not match" }

$idx = Find-ParaIndex("We have revised the paper in light")
if ($idx -lt 0) { Write-Output "WARN: could not locate 'We have revised the paper' paragraph" }
$d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
$d.Paragraphs.Item($idx + 1).Range.Text = "We enclose a revised version of the manuscript (clean and with tracked changes), which addresses the reviewers’ comments and suggestions. The reviewers’ comments are shown below with our response to each comment in italics."

$idx2 = Find-ParaIndex("We enclose a revised version of the manuscript (clean and with tracked changes), which addresses the reviewers")
if ($idx2 -lt 0) { Write-Output "WARN: could not locate re-inserted 'We enclose' paragraph" }
$d.Paragraphs.Item($idx2).Range.InsertParagraphAfter()
$d.Paragraphs.Item($idx2 + 1).Range.Text = "We hope that the manuscript is now suitable for publication in Vaccine."

# --- 3. Replace the "Whilst we agree that the use of ..." paragraph (italic,
#        multiple runs) with a new single-run italic sentence, then insert a
#        new italic paragraph after it with the remaining response text. ---
$ok3 = $d.Content.Find.Execute(
    "Whilst we agree that the use of “{% all cases}(% complete within vaccine status)[complete within category]” makes the table more complex to interpret we feel that this structure allows us to more clearly detail the missing data present. This is an important consideration as the volume of missing data is important to understand when interpreting the study findings.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We have also clarified the definition of sputum spear status in Table 1 and UK birth status in Table 2.",
    2)
if (-not $ok3) { Write-Output "WARN: step 3 find/replace did not match" }

$idx3 = Find-ParaIndex("We have also clarified the definition")
if ($idx3 -lt 0) { Write-Output "WARN: could not locate 'We have also clarified' paragraph" }
$d.Paragraphs.Item($idx3).Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($idx3 + 1)
$newPara.Range.Text = "We agree that the use of the {% all cases}(% complete within vaccine status)[complete within category] structure in Table 1 and Table 2 adds complexity but feel that this is justified as missing data is an important consideration for this study. We would welcome suggestions for improvements."

# --- 4. Reword the table-footer acknowledgement sentence ---
$ok4 = $d.Content.Find.Execute(
    "We agree that this was missing from the table. We have added the following definition to the table footer:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Thank you. We have added the following definition to the table footer:",
    2)
if (-not $ok4) { Write-Output "WARN: step 4 find/replace did not match" }

Write-Output "done"
